$d = $word.ActiveDocument

$oldText = "por Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/)."
$newText = "por Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

# Locate the whole "por Jenik Hollan, ..." sentence (it is split across many
# differently-formatted runs — plain text runs plus a hyperlink-styled run).
$rng = $d.Content.Duplicate
$rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Replace the whole run-fragmented sentence with a single, plainly
# formatted run containing the updated (2022) link text.
$rng.Delete()
$rng.Font.Reset()
$rng.InsertAfter($newText)
